$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.019492762824138
$ws.Range("D2").Value = 1.021393178369514
$ws.Range("E2").Value = 1.020624304830882
$ws.Range("F2").Value = 1.030775761595389
$ws.Range("I2").Value = 1.029441138464657
$ws.Range("J2").Value = 1.024695219552408
$ws.Range("K2").Value = 1.024230641327941
$ws.Range("L2").Value = 1.02346403984532
$ws.Range("M2").Value = 1.033585795494001
$ws.Range("N2").Value = 1.012311820580684
$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.020332698529388
$ws.Range("D3").Value = 1.022106664197465
$ws.Range("E3").Value = 1.021333856843096
$ws.Range("F3").Value = 1.031848139781905
$ws.Range("I3").Value = 1.029539744899054
$ws.Range("J3").Value = 1.02517250677231
$ws.Range("K3").Value = 1.02475082503839
$ws.Range("L3").Value = 1.023980136191531
$ws.Range("M3").Value = 1.034465890600347
$ws.Range("N3").Value = 1.012471120091427
$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.020876708629403
$ws.Range("D4").Value = 1.022569114425624
$ws.Range("E4").Value = 1.021793825751776
$ws.Range("F4").Value = 1.032542672403855
$ws.Range("I4").Value = 1.029602199551522
$ws.Range("J4").Value = 1.025481226512583
$ws.Range("K4").Value = 1.025087525795704
$ws.Range("L4").Value = 1.024314255369686
$ws.Range("M4").Value = 1.035035411886203
$ws.Range("N4").Value = 1.012574115757298
$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.021105532474409
$ws.Range("D5").Value = 1.0227637129182
$ws.Range("E5").Value = 1.021987396433201
$ws.Range("F5").Value = 1.032834804370152
$ws.Range("I5").Value = 1.029628131659726
$ws.Range("J5").Value = 1.025610983162145
$ws.Range("K5").Value = 1.025229099209609
$ws.Range("L5").Value = 1.024454758617546
$ws.Range("M5").Value = 1.035274847496945
$ws.Range("N5").Value = 1.012617395143297
$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.021143960117676
$ws.Range("D6").Value = 1.022796397627826
$ws.Range("E6").Value = 1.022019909467642
$ws.Range("F6").Value = 1.032883863409123
$ws.Range("I6").Value = 1.029632466769056
$ws.Range("J6").Value = 1.025632768133901
$ws.Range("K6").Value = 1.025252871404783
$ws.Range("L6").Value = 1.024478352014208
$ws.Range("M6").Value = 1.035315050282284
$ws.Range("N6").Value = 1.012624660758199
$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.020879765707223
$ws.Range("D7").Value = 1.022571713937644
$ws.Range("E7").Value = 1.021796411469894
$ws.Range("F7").Value = 1.032546575295442
$ws.Range("I7").Value = 1.029602547330743
$ws.Range("J7").Value = 1.025482960443148
$ws.Range("K7").Value = 1.025089417412197
$ws.Range("L7").Value = 1.024316132627124
$ws.Range("M7").Value = 1.035038611202894
$ws.Range("N7").Value = 1.01257469413765
$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.019776515428744
$ws.Range("D8").Value = 1.021634142336088
$ws.Range("E8").Value = 1.020863926059805
$ws.Range("F8").Value = 1.031138045372497
$ws.Range("I8").Value = 1.029474742222185
$ws.Range("J8").Value = 1.024856544558189
$ws.Range("K8").Value = 1.02440641693228
$ws.Range("L8").Value = 1.023638421013014
$ws.Range("M8").Value = 1.033883218726876
$ws.Range("N8").Value = 1.01236567327712
$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.017836458663807
$ws.Range("D9").Value = 1.019988046289175
$ws.Range("E9").Value = 1.019227284317939
$ws.Range("F9").Value = 1.028660917661355
$ws.Range("I9").Value = 1.029239216401485
$ws.Range("J9").Value = 1.023751872947269
$ws.Range("K9").Value = 1.023203757369157
$ws.Range("L9").Value = 1.022445565504003
$ws.Range("M9").Value = 1.031847623816531
$ws.Range("N9").Value = 1.011996743917428
$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.01654586848018
$ws.Range("D10").Value = 1.018894797412809
$ws.Range("E10").Value = 1.01814066992757
$ws.Range("F10").Value = 1.027012833600084
$ws.Range("I10").Value = 1.029075295992941
$ws.Range("J10").Value = 1.023014922111244
$ws.Range("K10").Value = 1.022402645143613
$ws.Range("L10").Value = 1.021651315319635
$ws.Range("M10").Value = 1.030490852594618
$ws.Range("N10").Value = 1.011750406958414
$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.015987705842455
$ws.Range("D11").Value = 1.018422413106312
$ws.Range("E11").Value = 1.017671237300681
$ws.Range("F11").Value = 1.026299996300818
$ws.Range("I11").Value = 1.02900268624977
$ws.Range("J11").Value = 1.022695708810437
$ws.Range("K11").Value = 1.022055926725105
$ws.Range("L11").Value = 1.021307645644763
$ws.Range("M11").Value = 1.029903435844896
$ws.Range("N11").Value = 1.011643654415615
$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.015780481651771
$ws.Range("D12").Value = 1.018247100390436
$ws.Range("E12").Value = 1.017497032844703
$ws.Range("F12").Value = 1.026035336845664
$ws.Range("I12").Value = 1.028975471235324
$ws.Range("J12").Value = 1.022577123546558
$ws.Range("K12").Value = 1.021927166598598
$ws.Range("L12").Value = 1.021180029599739
$ws.Range("M12").Value = 1.029685255065002
$ws.Range("N12").Value = 1.011603989107354
$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.015824927294046
$ws.Range("D13").Value = 1.018284698665032
$ws.Range("E13").Value = 1.017534392856111
$ws.Range("F13").Value = 1.026092101751943
$ws.Range("I13").Value = 1.028981320004885
$ws.Range("J13").Value = 1.022602561159728
$ws.Range("K13").Value = 1.021954784862593
$ws.Range("L13").Value = 1.021207401927393
$ws.Range("M13").Value = 1.029732055028749
$ws.Range("N13").Value = 1.011612498016937
$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.015970574531979
$ws.Range("D14").Value = 1.01840791859528
$ws.Range("E14").Value = 1.017656834156402
$ws.Range("F14").Value = 1.026278117014092
$ws.Range("I14").Value = 1.029000441633281
$ws.Range("J14").Value = 1.022685906819858
$ws.Range("K14").Value = 1.022045282824577
$ws.Range("L14").Value = 1.021297096076817
$ws.Range("M14").Value = 1.029885400704938
$ws.Range("N14").Value = 1.011640375924753
$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.016060326195001
$ws.Range("D15").Value = 1.018483858632054
$ws.Range("E15").Value = 1.017732296027673
$ws.Range("F15").Value = 1.026392743088237
$ws.Range("I15").Value = 1.029012190713993
$ws.Range("J15").Value = 1.02273725684657
$ws.Range("K15").Value = 1.02210104515802
$ws.Range("L15").Value = 1.021352364697793
$ws.Range("M15").Value = 1.029979883645826
$ws.Range("N15").Value = 1.011657550758497
$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.016582926161054
$ws.Range("D16").Value = 1.018926169207082
$ws.Range("E16").Value = 1.018171847492733
$ws.Range("F16").Value = 1.027060159057808
$ws.Range("I16").Value = 1.029080080542206
$ws.Range("J16").Value = 1.023036105050227
$ws.Range("K16").Value = 1.022425659356819
$ws.Range("L16").Value = 1.021674128829969
$ws.Range("M16").Value = 1.030529839141399
$ws.Range("N16").Value = 1.01175748997501
$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.016910920054279
$ws.Range("D17").Value = 1.019203887757633
$ws.Range("E17").Value = 1.018447856465593
$ws.Range("F17").Value = 1.027479024575189
$ws.Range("I17").Value = 1.029122229717063
$ws.Range("J17").Value = 1.023223536303456
$ws.Range("K17").Value = 1.02262932723283
$ws.Range("L17").Value = 1.021876029672101
$ws.Range("M17").Value = 1.030874832119586
$ws.Range("N17").Value = 1.011820156200225
$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.017102298032703
$ws.Range("D18").Value = 1.019365972632215
$ws.Range("E18").Value = 1.018608951784113
$ws.Range("F18").Value = 1.027723418299913
$ws.Range("I18").Value = 1.029146657201864
$ws.Range("J18").Value = 1.023332851217948
$ws.Range("K18").Value = 1.022748139431337
$ws.Range("L18").Value = 1.021993818660379
$ws.Range("M18").Value = 1.03107606778565
$ws.Range("N18").Value = 1.011856699925974
$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.017167563907919
$ws.Range("D19").Value = 1.019421255676668
$ws.Range("E19").Value = 1.018663898705973
$ws.Range("F19").Value = 1.027806763206029
$ws.Range("I19").Value = 1.029154959626197
$ws.Range("J19").Value = 1.023370122953334
$ws.Range("K19").Value = 1.022788654026526
$ws.Range("L19").Value = 1.02203398563422
$ws.Range("M19").Value = 1.031144685126663
$ws.Range("N19").Value = 1.011869158950812
$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.016875722712424
$ws.Range("D20").Value = 1.019174081208304
$ws.Range("E20").Value = 1.018418232550313
$ws.Range("F20").Value = 1.027434076345193
$ws.Range("I20").Value = 1.029117723787882
$ws.Range("J20").Value = 1.023203427788567
$ws.Range("K20").Value = 1.022607473915627
$ws.Range("L20").Value = 1.021854365168115
$ws.Range("M20").Value = 1.030837816897803
$ws.Range("N20").Value = 1.011813433575598
$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.015927682216371
$ws.Range("D21").Value = 1.018371629191926
$ws.Range("E21").Value = 1.017620773712202
$ws.Range("F21").Value = 1.026223336835322
$ws.Range("I21").Value = 1.02899481753604
$ws.Range("J21").Value = 1.022661364011875
$ws.Range("K21").Value = 1.022018632688159
$ws.Range("L21").Value = 1.021270682321709
$ws.Range("M21").Value = 1.029840243888728
$ws.Range("N21").Value = 1.01163216692193
$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.015332203893727
$ws.Range("D22").Value = 1.017867975127725
$ws.Range("E22").Value = 1.017120327770418
$ws.Range("F22").Value = 1.02546279184885
$ws.Range("I22").Value = 1.028916127028716
$ws.Range("J22").Value = 1.022320459425085
$ws.Range("K22").Value = 1.021648558946075
$ws.Range("L22").Value = 1.020903919221743
$ws.Range("M22").Value = 1.029213099723413
$ws.Range("N22").Value = 1.011518124334923
$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.015647821605052
$ws.Range("D23").Value = 1.018134887769414
$ws.Range("E23").Value = 1.017385533205783
$ws.Range("F23").Value = 1.025865904999657
$ws.Range("I23").Value = 1.028957976251964
$ws.Range("J23").Value = 1.022501187329482
$ws.Range("K23").Value = 1.021844727023283
$ws.Range("L23").Value = 1.021098325876821
$ws.Range("M23").Value = 1.02954555384244
$ws.Range("N23").Value = 1.011578587254647
$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.016891626676421
$ws.Range("D24").Value = 1.019187549207456
$ws.Range("E24").Value = 1.018431618001486
$ws.Range("F24").Value = 1.027454386279009
$ws.Range("I24").Value = 1.029119760309767
$ws.Range("J24").Value = 1.023212513993649
$ws.Range("K24").Value = 1.022617348439177
$ws.Range("L24").Value = 1.021864154352019
$ws.Range("M24").Value = 1.030854542461306
$ws.Range("N24").Value = 1.011816471266368
$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.01833752583675
$ws.Range("D25").Value = 1.020412877516959
$ws.Range("E25").Value = 1.019649613001743
$ws.Range("F25").Value = 1.029300730664936
$ws.Range("I25").Value = 1.029301324646869
$ws.Range("J25").Value = 1.024037550558034
$ws.Range("K25").Value = 1.023514562546928
$ws.Range("L25").Value = 1.022753778661254
$ws.Range("M25").Value = 1.03237382637433
$ws.Range("N25").Value = 1.012092190495369

Write-Host "Applied 264 value updates to vm_pu sheet (case with 380 kV)"
